# Auto-generated edit script: updates cryptos list values to match target snapshot
# (rankings/prices/volume-1h refreshed by the scheduled GitHub Actions scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number as text (e.g. "583.53") --
# format them as Text first so Excel keeps them as literal strings instead of
# silently converting to a numeric cell (matches source data's inline-string type).
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D11",
    "D12",
    "D14",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D30",
    "D31",
    "D33",
    "D35",
    "D36",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# New cell values (applies to every touched cell, text-forced or not)
$updates = [ordered]@{
    "D2" = "60.720.25"
    "E2" = "  +0.41%  "
    "D3" = "2.652.08"
    "E3" = "  +1.99%  "
    "E4" = "  +0.02%  "
    "D5" = "583.53"
    "E5" = "  +2.14%  "
    "D6" = "144.92"
    "E6" = "  +1.80%  "
    "D7" = "0.998"
    "E7" = "  +0.13%  "
    "E8" = "  -0.23%  "
    "D9" = "6.57"
    "E9" = "  +1.44%  "
    "E10" = "  +1.26%  "
    "D11" = "0.374"
    "E11" = "  +2.70%  "
    "D12" = "0.155"
    "E12" = "  +2.01%  "
    "D13" = "3.125.88"
    "E13" = "  +2.04%  "
    "D14" = "25.73"
    "E14" = "  +10.54%  "
    "D15" = "60.692.60"
    "E15" = "  +0.34%  "
    "E16" = "  +1.86%  "
    "D17" = "2.661.92"
    "D18" = "11.53"
    "E18" = "  +1.36%  "
    "D19" = "4.72"
    "E19" = "  +1.37%  "
    "D20" = "350.17"
    "E20" = "  +1.08%  "
    "D21" = "6.93"
    "E21" = "  -1.17%  "
    "D22" = "0.998"
    "E22" = "  +0.06%  "
    "D23" = "0.535"
    "E23" = "  +0.11%  "
    "D24" = "63.98"
    "E24" = "  +1.13%  "
    "D25" = "0.996"
    "E25" = "  +0.08%  "
    "D26" = "0.161"
    "E26" = "  +1.38%  "
    "D27" = "8.12"
    "E27" = "  +5.45%  "
    "D28" = "1.98"
    "E28" = "  +8.94%  "
    "D29" = "0.0₃0809"
    "E29" = "  +2.84%  "
    "D30" = "6.74"
    "E30" = "  +5.55%  "
    "D31" = "167.61"
    "E31" = "  +4.55%  "
    "E32" = "  +0.11%  "
    "D33" = "19.84"
    "E33" = "  +1.78%  "
    "E34" = "  +9.31%  "
    "D35" = "4.44"
    "E35" = "  +5.08%  "
    "D36" = "1.31"
    "E36" = "  +7.90%  "
    "E37" = "  +2.17%  "
    "D38" = "327.71"
    "E38" = "  +10.78%  "
    "E39" = "  +4.48%  "
    "D40" = "38.35"
    "E40" = "  +1.52%  "
    "D41" = "0.875"
    "E41" = "  +3.33%  "
    "D42" = "5.24"
    "E42" = "  +7.44%  "
    "D43" = "20.50"
    "E43" = "  +3.84%  "
    "D44" = "134.86"
    "E44" = "  -2.38%  "
    "B45" = "Stellar"
    "C45" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D45" = "0.0997"
    "E45" = "  +1.36%  "
    "B46" = "Mantle"
    "C46" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
    "D46" = "0.613"
    "E46" = "  +0.59%  "
    "B47" = "FirstDigitalUSD"
    "C47" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "D47" = "0.999"
    "E47" = "  +0.25%  "
    "B48" = "InjectiveProtocol"
    "C48" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D48" = "20.49"
    "E48" = "  +3.34%  "
    "B49" = "Hedera"
    "C49" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D49" = "0.0557"
    "E49" = "  +1.90%  "
    "D50" = "0.0245"
    "E50" = "  +2.20%  "
    "D51" = "2.126.46"
    "E51" = "  +4.98%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
